# NIT-9006960287.xlsx — "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# The sheet "Hoja1" lists worker mora (arrears) records. Column E ("Periodo
# Mora") holds the period code. The old period "2508" is retired and the
# three worker rows (16, 17, 18) are updated to the new period "2509". The
# period column for these rows is also center-aligned (matching the rest of
# the centered table columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Periodo Mora" value for each worker row: 2508 -> 2509
$ws.Range("E16").Value = "2509"
$ws.Range("E17").Value = "2509"
$ws.Range("E18").Value = "2509"

# Center the "Periodo Mora" column values for those rows
$ws.Range("E16:E18").HorizontalAlignment = -4108
